# Atualização de bases das ligas, do dia: 08-05-2024 às 20:15
#
# This script swaps the data of two pairs of rows in the "Israel Premier
# League" sheet (rows 219/221 and rows 222/223), which is the effect of
# re-sorting/re-ordering the source rows by match id/date during the
# league-base refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Row 219 <-> Row 221 : full row swap (id, teams, scores and all odds)
# ---------------------------------------------------------------------
$ws.Range("B219").Value = 8016156
$ws.Range("E219").Value = "Beitar Jerusalem"
$ws.Range("F219").Value = "MS Ashdod"
$ws.Range("G219").Value = 0
$ws.Range("H219").Value = 0
$ws.Range("J219").Value = 1.8
$ws.Range("K219").Value = 3.6
$ws.Range("L219").Value = 4
$ws.Range("M219").Value = 2
$ws.Range("N219").Value = 3.4
$ws.Range("O219").Value = 3.4
$ws.Range("P219").Value = -0.5
$ws.Range("Q219").Value = 2.05
$ws.Range("R219").Value = 1.8
$ws.Range("S219").Value = 2.5
$ws.Range("T219").Value = 2.05
$ws.Range("U219").Value = 1.8
$ws.Range("V219").Value = -1
$ws.Range("W219").Value = 2.4
$ws.Range("X219").Value = -1
$ws.Range("Y219").Value = -1
$ws.Range("Z219").Value = 0.8
$ws.Range("AA219").Value = -1
$ws.Range("AB219").Value = 0.8

$ws.Range("B221").Value = 8015672
$ws.Range("E221").Value = "Hapoel Bnei Sakhnin"
$ws.Range("F221").Value = "Maccabi Tel Aviv"
$ws.Range("G221").Value = 1
$ws.Range("H221").Value = 1
$ws.Range("J221").Value = 7
$ws.Range("K221").Value = 4.333
$ws.Range("L221").Value = 1.444
$ws.Range("M221").Value = 10
$ws.Range("N221").Value = 5.25
$ws.Range("O221").Value = 1.3
$ws.Range("P221").Value = 1.5
$ws.Range("Q221").Value = 1.975
$ws.Range("R221").Value = 1.875
$ws.Range("S221").Value = 3
$ws.Range("T221").Value = 1.975
$ws.Range("U221").Value = 1.875
$ws.Range("V221").Value = -1
$ws.Range("W221").Value = 4.25
$ws.Range("X221").Value = -1
$ws.Range("Y221").Value = 0.9750000000000001
$ws.Range("Z221").Value = -1
$ws.Range("AA221").Value = -1
$ws.Range("AB221").Value = 0.875

# ---------------------------------------------------------------------
# Row 222 <-> Row 223 : ids 8015674 / 8015675 swap (shared-string ids),
# teams swap, and odds refresh
# ---------------------------------------------------------------------
# these "id" values are stored as text (not numbers) in the sheet, so
# force a text format before writing the digit-only string
$ws.Range("B222").NumberFormat = "@"
$ws.Range("B222").Value = "8015675"
$ws.Range("E222").Value = "Hapoel Beer Sheva"
$ws.Range("F222").Value = "Maccabi Bnei Raina"
$ws.Range("J222").Value = 1.65
$ws.Range("K222").Value = 3.6
$ws.Range("L222").Value = 5.25
$ws.Range("M222").Value = 1.5
$ws.Range("N222").Value = 4
$ws.Range("O222").Value = 6.5
$ws.Range("P222").Value = -1
$ws.Range("Q222").Value = 1.875
$ws.Range("R222").Value = 1.975
$ws.Range("S222").Value = 2.25
$ws.Range("T222").Value = 1.8
$ws.Range("U222").Value = 2.05

$ws.Range("B223").NumberFormat = "@"
$ws.Range("B223").Value = "8015674"
$ws.Range("E223").Value = "Hapoel Haifa"
$ws.Range("F223").Value = "Hapoel Bnei Sakhnin"
$ws.Range("J223").Value = 1.95
$ws.Range("K223").Value = 3.4
$ws.Range("L223").Value = 3.75
$ws.Range("M223").Value = 1.95
$ws.Range("N223").Value = 3.4
$ws.Range("O223").Value = 3.75
$ws.Range("P223").Value = -0.5
$ws.Range("Q223").Value = 2.025
$ws.Range("R223").Value = 1.825
$ws.Range("S223").Value = 2.25
$ws.Range("T223").Value = 1.85
$ws.Range("U223").Value = 2
